$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New column D ("Chú thích" / remarks) header, styled like the other
#    header cells (A1:C1).
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "Chú thích"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Columns.Item(4).ColumnWidth = 83

# ---------------------------------------------------------------------------
# 2) Normalise borders on the existing data rows (2-6): date (A) and name
#    (B) columns pick up the same bordered/centred look already used on row 2.
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A3:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B2").Copy()
$ws.Range("B3:B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Totals column (C): thousands separator + border + centred, applied to the
# whole existing column of values first …
$ws.Range("C2").Borders.LineStyle = 1
$ws.Range("C2").NumberFormat = "#,##0"
$ws.Range("C2").HorizontalAlignment = -4108   # xlCenter
$ws.Range("C2").VerticalAlignment = -4108     # xlCenter
$ws.Range("C2").Copy()
$ws.Range("C3:C6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New notes column (D): bordered, blank for the pre-existing rows.
$ws.Range("D2").Borders.LineStyle = 1
$ws.Range("D2").Copy()
$ws.Range("D3:D6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) New row 7 — MicroSD card purchased for the Raspberry Pi project.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 43187          # 2018-03-28, serial date value
$ws.Range("B7").Value = "MicroSD 16GB"
$ws.Range("C7").Value = 170000
$ws.Range("D7").Value = "https://memoryzone.com.vn/the-nho-microsdhc-sandisk-ultra-16gb-80mbs-533x-2017/"

$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B2").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C2").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D2").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4) Leave the selection cursor where the source workbook had it.
# ---------------------------------------------------------------------------
[void]$ws.Range("F5").Select()
